$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Remove the "Close Percentages" column N (header comment, header cell, and data) ---
$ws.Range("N1").Comment.Delete()

# Re-style N1 to match the other (unused) header cells instead of the special bordered/wrapped style
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").ClearContents()

# Clear the "Close Percentages" data cells in N2:N4 entirely (value + formatting)
$ws.Range("N2:N4").Clear()

# The header row no longer needs the extra height used for the wrapped "Close Percentages" text
$ws.Rows.Item(1).RowHeight = 12.75

# --- 2. Add a new data row (Capital Call 9) below the existing rows ---
$ws.Range("A4:M4").Copy()
$ws.Range("A5:M5").PasteSpecial(-4122)

$ws.Range("A5").Value = "SAAS Fund"
$ws.Range("B5").Value = "CoInvest"
$ws.Range("C5").Value = "Capital Call 9"
$ws.Range("D5").Value = "30, 10"
$ws.Range("E5").Value = 44905
$ws.Range("F5").Value = 44816
$ws.Range("G5").Value = "Yes"
$ws.Range("H5").Value = "No"
$ws.Range("I5").Value = "First Close, Second Close"
$ws.Range("J5").Value = "Percentage of Commitment"
$ws.Range("K5").Value = "Series A:1000:0,Series B:1100:100,Series C:1200:200"
$ws.Range("L5").Value = "Yes"
$ws.Range("M5").Value = "Yes"
$ws.Rows.Item(5).RowHeight = 13.8

# --- 3. Update the view: scroll so column F is the leftmost visible column and select N5 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 6
$ws.Range("N5").Select()

Write-Host "Edit complete"
